# The "Förändrad" (Changed) date in column C was bumped from 2023-09-21
# (serial 45190) to 2023-09-23 (serial 45192) for every data row
# (rows 2 through 307) on the single worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 307
$col = 3  # column C

$ws.Range($ws.Cells.Item($firstRow, $col), $ws.Cells.Item($lastRow, $col)).Value = 45192
